# Generate Report for Handoff
# Refresh the localization-status report: the "113fd6e5-d39d-46be-b761-ffd22f2f26ea"
# file has been re-handed-off, so its handoff timestamp moves forward and its
# priority (along with the other remaining files' priority) is recalculated
# from "low" to "ht".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: the "Latest HO Xliff Generate Date" shown for the
# "Ready for handoff" rows (113fd6e5, 4010d2ce, 819faaa8, bc524db0 - rows 4-7)
# all shared the same timestamp value, which moves forward.
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-04 14:34:47"
}

# zh-cn sheet: rows 4-7 (113fd6e5, 4010d2ce, 819faaa8, bc524db0) get a new
# Priority value and an updated Latest Handoff Datetime.
foreach ($r in 4..7) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-04 14:34:43"
}

# de-de sheet: rows 4-7 get the same new Priority value, and the Latest
# Handoff Datetime (which shared the same underlying timestamp as the
# Overview sheet) also moves forward to the same new value.
foreach ($r in 4..7) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-04 14:34:47"
}
